$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.72%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.81%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.165"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.11%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07478"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.82%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.768"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.00%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.668"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.17%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.800"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.46%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9250"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.44%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1710"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.29%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07502"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.45%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07947"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.41%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.94%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09887"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.87%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04666"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.67%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006241"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.17%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.96%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.49%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.61%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1347"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.99%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.563"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.47%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1551"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.37%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.21%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004414"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.58%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001401"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "19.92%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001810"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "8.65%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01651"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04515"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.34%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006940"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-4.59%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1341"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.88%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002061"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.63%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01326"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.17%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006077"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.79%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.930"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.98%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01225"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-5.79%"
